$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "26.986.72"
$ws.Range("E2").Value = "  +2.03%  "
$ws.Range("D3").Value = "1.815.17"
$ws.Range("E3").Value = "  +2.52%  "
$ws.Range("E4").Value = "  +0.23%  "
Set-TextValue $ws.Range("D5") "312.96"
$ws.Range("E5").Value = "  +2.04%  "
Set-TextValue $ws.Range("D6") "1.006"
Set-TextValue $ws.Range("D7") "0.4292"
$ws.Range("E7").Value = "  -0.08%  "
Set-TextValue $ws.Range("D8") "0.3668"
$ws.Range("E8").Value = "  +0.09%  "
Set-TextValue $ws.Range("D9") "0.07254"
$ws.Range("E9").Value = "  +0.67%  "
$ws.Range("D10").Value = "2.188.06"
$ws.Range("E10").Value = "  +24.41%  "
Set-TextValue $ws.Range("D11") "0.8628"
$ws.Range("E11").Value = "  +1.61%  "
$ws.Range("E12").Value = "  +4.56%  "
Set-TextValue $ws.Range("D13") "5.408"
$ws.Range("E13").Value = "  +3.33%  "
Set-TextValue $ws.Range("D14") "6.600"
$ws.Range("E14").Value = "  +2.64%  "
Set-TextValue $ws.Range("D15") "0.06941"
$ws.Range("E15").Value = "  +0.66%  "
Set-TextValue $ws.Range("D16") "81.10"
$ws.Range("E16").Value = "  +2.22%  "
Set-TextValue $ws.Range("D17") "1.012"
$ws.Range("E17").Value = "  +0.60%  "
Set-TextValue $ws.Range("D18") "0.000008892"
$ws.Range("E18").Value = "  +2.39%  "
$ws.Range("E19").Value = "  +0.37%  "
Set-TextValue $ws.Range("D20") "15.18"
$ws.Range("E20").Value = "  +0.87%  "
$ws.Range("D21").Value = "27.021.49"
$ws.Range("E21").Value = "  +2.18%  "
Set-TextValue $ws.Range("D22") "5.172"
$ws.Range("E22").Value = "  +1.41%  "
$ws.Range("D23").Value = "2.425.36"
$ws.Range("E23").Value = "  +22.68%  "
Set-TextValue $ws.Range("D24") "10.99"
$ws.Range("E24").Value = "  -2.51%  "
Set-TextValue $ws.Range("D25") "153.83"
$ws.Range("E25").Value = "  +1.06%  "
$ws.Range("E26").Value = "  +0.04%  "
Set-TextValue $ws.Range("D27") "18.33"
$ws.Range("E27").Value = "  +1.02%  "
Set-TextValue $ws.Range("D28") "5.221"
$ws.Range("E28").Value = "  +2.49%  "
Set-TextValue $ws.Range("D29") "1.900"
$ws.Range("E29").Value = "  +9.73%  "
Set-TextValue $ws.Range("D30") "114.68"
$ws.Range("E30").Value = "  +0.12%  "
Set-TextValue $ws.Range("D31") "0.08933"
$ws.Range("E31").Value = "  -0.40%  "
Set-TextValue $ws.Range("D32") "1.187"
$ws.Range("E32").Value = "  +6.83%  "
Set-TextValue $ws.Range("D33") "0.7468"
$ws.Range("E33").Value = "  +3.30%  "
Set-TextValue $ws.Range("D34") "4.415"
$ws.Range("E34").Value = "  +2.20%  "
Set-TextValue $ws.Range("D35") "2.812"
$ws.Range("E35").Value = "  +2.35%  "
Set-TextValue $ws.Range("D36") "1.006"
$ws.Range("E36").Value = "  +0.34%  "
Set-TextValue $ws.Range("D37") "1.132"
$ws.Range("E37").Value = "  +4.88%  "
$ws.Range("E38").Value = "  +0.77%  "
Set-TextValue $ws.Range("D39") "0.01921"
$ws.Range("E39").Value = "  +1.55%  "
Set-TextValue $ws.Range("D40") "0.5099"
$ws.Range("E40").Value = "  +3.63%  "
Set-TextValue $ws.Range("D41") "2.744"
$ws.Range("E41").Value = "  +6.58%  "
$ws.Range("E42").Value = "  +3.11%  "
Set-TextValue $ws.Range("D43") "6.473"
$ws.Range("E43").Value = "  +3.56%  "
Set-TextValue $ws.Range("D44") "8.306"
$ws.Range("E44").Value = "  +3.67%  "
Set-TextValue $ws.Range("D47") "1.006"
$ws.Range("E47").Value = "  +0.36%  "
Set-TextValue $ws.Range("D48") "0.4581"
$ws.Range("E48").Value = "  +2.27%  "
Set-TextValue $ws.Range("D49") "1.645"
$ws.Range("E49").Value = "  +3.46%  "
Set-TextValue $ws.Range("D50") "0.06211"
$ws.Range("E50").Value = "  +0.20%  "
Set-TextValue $ws.Range("D51") "1.841"
$ws.Range("E51").Value = "  +5.96%  "

# Row 45/46 swap: EnergySwap <-> Quant
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D45") "106.61"
$ws.Range("E45").Value = "  +1.75%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D46") "10.41"
$ws.Range("E46").Value = "  +3.03%  "
